# Apply the "Updated symbol list" crypto price refresh described by the diff.
# The source workbook stores every data cell as literal TEXT (inline strings),
# even numeric-looking values like "326.65" or "-2.38%". A leading apostrophe
# forces Excel's COM layer to keep the new value as text instead of inferring
# a number, and resetting the cell style back to "Normal" afterwards avoids
# leaving behind a "quote prefix" / Text-number-format style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Value)
    $Cell.Value = "'" + $Value
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "327.36"
Set-TextValue $ws.Range("E2") "-2.01%"
Set-TextValue $ws.Range("D3") "44.25"
Set-TextValue $ws.Range("E3") "0.69%"
Set-TextValue $ws.Range("D4") "5.561"
Set-TextValue $ws.Range("E4") "-3.12%"
Set-TextValue $ws.Range("D5") "0.08049"
Set-TextValue $ws.Range("E5") "-3.70%"
Set-TextValue $ws.Range("B6") "FTXToken"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D6") "1.896"
Set-TextValue $ws.Range("E6") "-2.53%"
Set-TextValue $ws.Range("B7") "GateToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.288"
Set-TextValue $ws.Range("E7") "-5.18%"
Set-TextValue $ws.Range("D9") "0.9445"
Set-TextValue $ws.Range("E9") "-0.18%"
Set-TextValue $ws.Range("D10") "0.1160"
Set-TextValue $ws.Range("E10") "-6.18%"
Set-TextValue $ws.Range("D11") "0.1834"
Set-TextValue $ws.Range("E11") "-6.86%"
Set-TextValue $ws.Range("D12") "0.09685"
Set-TextValue $ws.Range("E12") "-2.78%"
Set-TextValue $ws.Range("D13") "0.04359"
Set-TextValue $ws.Range("E13") "-0.93%"
Set-TextValue $ws.Range("E14") "-0.45%"
Set-TextValue $ws.Range("D15") "0.001282"
Set-TextValue $ws.Range("E15") "-1.23%"
Set-TextValue $ws.Range("B16") "CoinExToken"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D16") "0.04214"
Set-TextValue $ws.Range("E16") "-4.66%"
Set-TextValue $ws.Range("B17") "TigerCash"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.005994"
Set-TextValue $ws.Range("E17") "-1.56%"
Set-TextValue $ws.Range("B18") "LEO"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.620"
Set-TextValue $ws.Range("E18") "3.87%"
Set-TextValue $ws.Range("B19") "BitpandaEcosystemToken"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Range("D19") "0.3496"
Set-TextValue $ws.Range("E19") "-1.16%"
Set-TextValue $ws.Range("B20") "MCDex"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D20") "8.571"
Set-TextValue $ws.Range("E20") "-1.77%"
Set-TextValue $ws.Range("B21") "ProBitToken"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Range("D21") "0.1378"
Set-TextValue $ws.Range("E21") "0.97%"
Set-TextValue $ws.Range("B22") "ZBToken"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws.Range("D22") "0.2524"
Set-TextValue $ws.Range("E22") "-3.59%"
Set-TextValue $ws.Range("D23") "0.001248"
Set-TextValue $ws.Range("E23") "0.16%"
Set-TextValue $ws.Range("D24") "0.004486"
Set-TextValue $ws.Range("E24") "2.84%"
Set-TextValue $ws.Range("D25") "0.0001260"
Set-TextValue $ws.Range("E25") "-0.17%"
Set-TextValue $ws.Range("D26") "0.0003989"
Set-TextValue $ws.Range("E26") "-0.34%"
Set-TextValue $ws.Range("D38") "0.02625"
Set-TextValue $ws.Range("E38") "-6.46%"
Set-TextValue $ws.Range("D39") "0.05446"
Set-TextValue $ws.Range("E39") "-7.14%"
Set-TextValue $ws.Range("D40") "0.007587"
Set-TextValue $ws.Range("E40") "-4.11%"
Set-TextValue $ws.Range("D42") "0.007274"
Set-TextValue $ws.Range("E42") "-19.59%"
Set-TextValue $ws.Range("D43") "0.002016"
Set-TextValue $ws.Range("E43") "-6.04%"
Set-TextValue $ws.Range("D44") "0.008886"
Set-TextValue $ws.Range("E44") "-14.13%"
Set-TextValue $ws.Range("E45") "-4.43%"
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("E46") "-0.21%"
Set-TextValue $ws.Range("D47") "0.003553"
Set-TextValue $ws.Range("E47") "11.29%"
Set-TextValue $ws.Range("D48") "0.002270"
Set-TextValue $ws.Range("E48") "-0.34%"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("E49") "-0.21%"
Set-TextValue $ws.Range("D50") "0.0002000"
Set-TextValue $ws.Range("E50") "-0.21%"
